# Regenerate merged AHB files
# - Rename header labels from the "_old" / "_new" suffix convention to the
#   "_FV2410" / "_FV2504" suffix convention used by the newer file-version pair.
# - Turn the data range into a real Excel Table (Table1).
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row labels -------------------------------------
$headerRng = $ws.Range("A1:U1")
$headerRng.Replace("_old", "_FV2410", 2, 1, $false, $false, $false, $false)
$headerRng.Replace("_new", "_FV2504", 2, 1, $false, $false, $false, $false)

# --- 2. Convert the used range into an Excel Table ------------------------
$lastCell = $ws.Cells.SpecialCells(11)
$lastRow = $lastCell.Row
$lastCol = $lastCell.Column
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
